$d = $word.ActiveDocument

$replacements = @(
    @{old="76×65="; new="58×78="},
    @{old="34×38="; new="58×81="},
    @{old="62×88="; new="81×25="},
    @{old="89×38="; new="29×70="},
    @{old="89×52="; new="85×23="},
    @{old="56×49="; new="79×99="},
    @{old="66×19="; new="69×63="},
    @{old="22×82="; new="90×13="},
    @{old="92×26="; new="32×31="},
    @{old="35×43="; new="44×70="},
    @{old="97×44="; new="37×33="},
    @{old="36×43="; new="93×19="},
    @{old="67×45="; new="39×20="},
    @{old="39×46="; new="42×46="},
    @{old="15×68="; new="90×14="},
    @{old="72×94="; new="92×51="},
    @{old="46×66="; new="97×23="},
    @{old="75×42="; new="55×85="},
    @{old="29×29="; new="47×60="},
    @{old="11×18="; new="82×68="},
    @{old="51×38="; new="21×95="},
    @{old="68×14="; new="74×12="},
    @{old="75×72="; new="18×90="},
    @{old="50×85="; new="62×47="},
    @{old="81×40="; new="28×86="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
